# Commit: "Updated status on 24/2/2021 in video/status_report"
#
# The weekly status-report sheet (Sheet1) gets two new rows appended:
#   row 11 - status as of 19/2/2021 (today's work becomes next row's "yesterday")
#   row 12 - status as of 24/2/2021
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 : 19/2/2021 status ---
$ws.Range("A11").Value = "19/2/2021"
$ws.Range("B11").Value = "Studied Basic tutorial 1 in gstreamers and `nAPI’s used in basic_tutorial_1 "
$ws.Range("C11").Value = "studied Basic tutorial 2 in gstreamers and `n API’s used in basic_tutorial_2"
$ws.Range("D11").Value = "Need to check playback tutorials"

# --- Row 12 : 24/2/2021 status ---
$ws.Range("A12").Value = "24/2/2021"
$ws.Range("B12").Value = "studied Basic tutorial 2 in gstreamers and `n API’s used in basic_tutorial_2"
$ws.Range("C12").Value = "studied playback turorial-1 and switching `nbetween audio streams"
$ws.Range("D12").Value = "Need to add logs to playback tutorial_1 and `nCheck for more info."

# --- Row heights, matching the taller wrapped text in these rows ---
$ws.Rows.Item(11).RowHeight = 48
$ws.Rows.Item(12).RowHeight = 54.75

# --- Column A ("Date") keeps the sheet's existing centered text-date look ---
$rngA = $ws.Range("A11:A12")
$rngA.HorizontalAlignment = -4108
$rngA.VerticalAlignment = -4108
$rngA.NumberFormat = "@"

# --- Column B ("Yesterday") : centered, wrapped, Calibri 11 ---
$rngB = $ws.Range("B11:B12")
$rngB.Font.Name = "Calibri"
$rngB.Font.Size = 11
$rngB.HorizontalAlignment = -4108
$rngB.VerticalAlignment = -4108
$rngB.WrapText = $true

# --- Column C ("Today") : centered, wrapped (same look used elsewhere in the sheet) ---
$rngC = $ws.Range("C11:C12")
$rngC.HorizontalAlignment = -4108
$rngC.VerticalAlignment = -4108
$rngC.WrapText = $true

# --- Column D ("Tomorrow") : D11 plain (like D10), D12 centered + wrapped ---
$rngD11 = $ws.Range("D11")
$rngD11.HorizontalAlignment = -4108
$rngD11.VerticalAlignment = -4108
$rngD11.WrapText = $false

$rngD12 = $ws.Range("D12")
$rngD12.HorizontalAlignment = -4108
$rngD12.VerticalAlignment = -4108
$rngD12.WrapText = $true

# --- Leave the selection on the last cell the author touched, D12 ---
$null = $ws.Range("D12").Select()
